# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q1" and "总计" holding the
#    per-fund breakdown for the new quarter (mirrors the layout already used
#    on the "2021-Q1" sheet).
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (2022-Q1: 3 funds, 2.42亿元), pushing the existing "2021-Q1" summary row
#    down by one.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet right after "2021-Q1"
# ------------------------------------------------------------------
$ws2021 = $wb.Worksheets.Item("2021-Q1")
$ws2022 = $wb.Worksheets.Add($null, $ws2021)
$ws2022.Name = "2022-Q1"

# Header row (B1:H1)
$headerRange = $ws2022.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

# Index column (A2:A4) carries the same header-like style as the source sheet
$indexRange = $ws2022.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Columns D:G hold figures that are stored as plain text (mirrors "2021-Q1")
$ws2022.Range("D2:G4").NumberFormat = "@"

# code, name, scale, totalPosition, positionShare, marketValue, rank
$fundRows = @(
    @("001481", "华宝油气(QDII)美元",    "39.80", "94.60", "2.61", "1.0388", 1),
    @("162411", "华宝油气(QDII)人民币A", "39.80", "94.60", "2.61", "1.0388", 1),
    @("007844", "华宝油气(QDII)人民币C", "12.98", "94.60", "2.61", "0.3388", 1)
)

for ($r = 0; $r -lt $fundRows.Length; $r++) {
    $data = $fundRows[$r]
    $excelRow = 2 + $r

    $ws2022.Cells.Item($excelRow, 1).Value = $r
    $ws2022.Cells.Item($excelRow, 2).NumberFormat = "@"
    $ws2022.Cells.Item($excelRow, 2).Value = $data[0]
    $ws2022.Cells.Item($excelRow, 3).Value = $data[1]
    $ws2022.Cells.Item($excelRow, 4).Value = $data[2]
    $ws2022.Cells.Item($excelRow, 5).Value = $data[3]
    $ws2022.Cells.Item($excelRow, 6).Value = $data[4]
    $ws2022.Cells.Item($excelRow, 7).Value = $data[5]
    $ws2022.Cells.Item($excelRow, 8).Value = $data[6]
}

# ------------------------------------------------------------------
# Step 2: insert the new 2022-Q1 summary row into "总计"
#         (re-fetch the sheet reference now that the sheet collection
#          has shifted because of the insert above)
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

$a2 = $wsTotal.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 3
$wsTotal.Cells.Item(2, 4).Value = 2.42

# renumber the (now pushed-down) original row so the index column stays 0-based
$wsTotal.Cells.Item(3, 1).Value = 1

# restore the originally-active sheet/selection
$ws2021.Activate()
$null = $ws2021.Range("A1").Select()

